$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '58.365.87'
Set-TextCell 2 5 '  -1.35%  '

Set-TextCell 3 4 '2.478.60'
Set-TextCell 3 5 '  -1.74%  '

Set-TextCell 4 4 '0.998'
Set-TextCell 4 5 '  +0.06%  '

Set-TextCell 5 4 '522.50'
Set-TextCell 5 5 '  -2.58%  '

Set-TextCell 6 4 '132.80'
Set-TextCell 6 5 '  -3.71%  '

Set-TextCell 7 4 '0.998'
Set-TextCell 7 5 '  +0.00%  '

Set-TextCell 8 4 '0.560'
Set-TextCell 8 5 '  -1.33%  '

Set-TextCell 9 4 '0.0997'
Set-TextCell 9 5 '  -1.57%  '

Set-TextCell 10 4 '0.157'
Set-TextCell 10 5 '  -0.79%  '

Set-TextCell 11 4 '5.40'
Set-TextCell 11 5 '  +0.79%  '

Set-TextCell 12 4 '0.344'
Set-TextCell 12 5 '  -1.34%  '

Set-TextCell 13 4 '2.913.25'
Set-TextCell 13 5 '  -1.47%  '

Set-TextCell 14 4 '58.260.37'
Set-TextCell 14 5 '  -1.31%  '

Set-TextCell 15 4 '22.19'
Set-TextCell 15 5 '  -3.65%  '

Set-TextCell 16 4 '0.0000137'
Set-TextCell 16 5 '  -1.75%  '

Set-TextCell 17 4 '2.476.26'
Set-TextCell 17 5 '  -1.90%  '

Set-TextCell 18 4 '10.89'
Set-TextCell 18 5 '  -2.10%  '

Set-TextCell 19 4 '321.72'
Set-TextCell 19 5 '  -1.12%  '

Set-TextCell 20 4 '4.19'
Set-TextCell 20 5 '  -2.30%  '

Set-TextCell 21 5 '  -0.14%  '

Set-TextCell 22 4 '5.79'
Set-TextCell 22 5 '  -2.95%  '

Set-TextCell 23 4 '64.56'
Set-TextCell 23 5 '  -1.87%  '

Set-TextCell 24 4 '0.410'
Set-TextCell 24 5 '  -3.12%  '

Set-TextCell 25 4 '1.00'
Set-TextCell 25 5 '  +0.06%  '

Set-TextCell 26 4 '0.162'
Set-TextCell 26 5 '  -3.32%  '

Set-TextCell 27 4 '7.43'
Set-TextCell 27 5 '  -2.99%  '

Set-TextCell 28 4 '0.0₃0753'
Set-TextCell 28 5 '  -2.27%  '

Set-TextCell 29 4 '6.39'
Set-TextCell 29 5 '  -4.37%  '

Set-TextCell 30 2 'PancakeSwap'
Set-TextCell 30 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 30 4 '1.71'
Set-TextCell 30 5 '  -4.23%  '

Set-TextCell 31 2 'Monero'
Set-TextCell 31 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 31 4 '167.37'
Set-TextCell 31 5 '  +2.14%  '

Set-TextCell 32 4 '1.17'
Set-TextCell 32 5 '  -0.93%  '

Set-TextCell 33 4 '0.998'

Set-TextCell 34 5 '  -0.15%  '

Set-TextCell 35 4 '18.20'
Set-TextCell 35 5 '  -1.53%  '

Set-TextCell 36 4 '1.34'
Set-TextCell 36 5 '  -9.42%  '

Set-TextCell 37 4 '4.01'
Set-TextCell 37 5 '  -2.57%  '

Set-TextCell 38 4 '1.50'
Set-TextCell 38 5 '  -3.54%  '

Set-TextCell 39 4 '0.799'
Set-TextCell 39 5 '  -2.60%  '

Set-TextCell 40 2 'RenderToken'
Set-TextCell 40 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 40 4 '5.15'
Set-TextCell 40 5 '  -1.11%  '

Set-TextCell 41 4 '277.54'
Set-TextCell 41 5 '  -3.11%  '

Set-TextCell 42 2 'Filecoin'
Set-TextCell 42 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 42 4 '3.49'
Set-TextCell 42 5 '  -4.18%  '

Set-TextCell 43 4 '0.599'
Set-TextCell 43 5 '  -1.05%  '

Set-TextCell 44 4 '126.55'
Set-TextCell 44 5 '  -4.70%  '

Set-TextCell 45 4 '0.0909'
Set-TextCell 45 5 '  -2.33%  '

Set-TextCell 46 4 '0.0494'
Set-TextCell 46 5 '  -3.16%  '

Set-TextCell 47 4 '0.0216'
Set-TextCell 47 5 '  -2.61%  '

Set-TextCell 48 4 '17.23'
Set-TextCell 48 5 '  -0.70%  '

Set-TextCell 49 4 '1.744.70'
Set-TextCell 49 5 '  -1.27%  '

Set-TextCell 50 4 '0.974'
Set-TextCell 50 5 '  -1.70%  '

Set-TextCell 51 5 '  -1.56%  '
